# Updates cryptos list values (price/volume columns, and a few row
# re-orderings / coin swaps) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.090.93"
$ws.Range("E2").Value = "  +5.28%  "
# Row 3
$ws.Range("D3").Value = "3.549.64"
$ws.Range("E3").Value = "  +4.50%  "
# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
# Row 5
$ws.Range("D5").Value = "'589.57"
$ws.Range("E5").Value = "  +5.17%  "
# Row 6
$ws.Range("D6").Value = "'192.20"
$ws.Range("E6").Value = "  +8.92%  "
# Row 7
$ws.Range("D7").Value = "'0.640"
$ws.Range("E7").Value = "  +1.25%  "
# Row 8
$ws.Range("D8").Value = "3.550.80"
$ws.Range("E8").Value = "  +4.81%  "
# Row 9
$ws.Range("E9").Value = "  -0.10%  "
# Row 10
$ws.Range("D10").Value = "'0.179"
$ws.Range("E10").Value = "  +3.39%  "
# Row 11
$ws.Range("E11").Value = "  +2.86%  "
# Row 12
$ws.Range("D12").Value = "'58.26"
$ws.Range("E12").Value = "  +8.36%  "
# Row 13
$ws.Range("E13").Value = "  +4.93%  "
# Row 14
$ws.Range("E14").Value = "  +4.47%  "
# Row 15
$ws.Range("D15").Value = "4.104.58"
$ws.Range("E15").Value = "  +4.10%  "
# Row 16
$ws.Range("E16").Value = "  +4.61%  "
# Row 17
$ws.Range("D17").Value = "3.550.48"
$ws.Range("E17").Value = "  +4.83%  "
# Row 18
$ws.Range("D18").Value = "69.112.95"
$ws.Range("E18").Value = "  +5.47%  "
# Row 19
$ws.Range("D19").Value = "'12.39"
$ws.Range("E19").Value = "  +4.35%  "
# Row 20
$ws.Range("D20").Value = "'0.120"
$ws.Range("E20").Value = "  +0.36%  "
# Row 21
$ws.Range("E21").Value = "  +3.39%  "
# Row 22
$ws.Range("E22").Value = "  +2.57%  "
# Row 23
$ws.Range("D23").Value = "'5.60"
$ws.Range("E23").Value = "  +12.95%  "
# Row 24
$ws.Range("E24").Value = "  +21.34%  "
# Row 25
$ws.Range("D25").Value = "'4.45"
$ws.Range("E25").Value = "  +8.01%  "
# Row 26
$ws.Range("D26").Value = "'90.98"
$ws.Range("E26").Value = "  +1.58%  "
# Row 27
$ws.Range("E27").Value = "  +4.05%  "
# Row 28
$ws.Range("D28").Value = "'11.14"
$ws.Range("E28").Value = "  +4.29%  "
# Row 29
$ws.Range("E29").Value = "  +5.80%  "
# Row 30
$ws.Range("D30").Value = "'31.85"
$ws.Range("E30").Value = "  +1.29%  "
# Row 31
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  +13.12%  "
# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'12.04"
$ws.Range("E32").Value = "  +4.39%  "
# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "'610.36"
$ws.Range("E33").Value = "  +5.71%  "
# Row 34
$ws.Range("D34").Value = "'65.33"
$ws.Range("E34").Value = "  +4.50%  "
# Row 35
$ws.Range("E35").Value = "  +5.91%  "
# Row 36
$ws.Range("D36").Value = "'0.148"
$ws.Range("E36").Value = "  +4.96%  "
# Row 37
$ws.Range("E37").Value = "  +0.07%  "
# Row 38
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0811"
$ws.Range("E38").Value = "  +9.51%  "
# Row 39
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").Value = "'37.60"
$ws.Range("E39").Value = "  +4.62%  "
# Row 40
$ws.Range("E40").Value = "  +5.37%  "
# Row 41
$ws.Range("D41").Value = "'3.56"
$ws.Range("E41").Value = "  -1.24%  "
# Row 42
$ws.Range("D42").Value = "3.286.08"
$ws.Range("E42").Value = "  +5.70%  "
# Row 43
$ws.Range("D43").Value = "'3.06"
$ws.Range("E43").Value = "  +9.30%  "
# Row 44
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "  +9.71%  "
# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0440"
$ws.Range("E45").Value = "  +5.13%  "
# Row 46
$ws.Range("D46").Value = "'3.27"
$ws.Range("E46").Value = "  +3.38%  "
# Row 47
$ws.Range("E47").Value = "  +1.77%  "
# Row 48
$ws.Range("D48").Value = "'2.76"
$ws.Range("E48").Value = "  +18.46%  "
# Row 49
$ws.Range("D49").Value = "'9.02"
$ws.Range("E49").Value = "  +6.58%  "
# Row 50
$ws.Range("E50").Value = "  +0.10%  "
# Row 51
$ws.Range("B51").Value = "LidoDAOToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D51").Value = "'3.21"
$ws.Range("E51").Value = "  +4.87%  "

# The apostrophe prefix above also sets a "quote prefix" flag on the
# cell style (the little green corner marker in Excel); resetting the
# style back to Normal clears that so only the cell VALUE changes,
# matching the source data (which carries no style changes here).
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
